$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.073399999999997
$ws.Range("A4").Value = -21.0046
$ws.Range("B4").Value = 5.437400000000004
$ws.Range("D4").Value = -6.677299999999994
$ws.Range("B5").Value = 5.413199999999996
$ws.Range("A6").Value = -21.2143
$ws.Range("A7").Value = -21.6638
$ws.Range("B8").Value = 4.964300000000001
$ws.Range("D9").Value = -8.3818
$ws.Range("D11").Value = -8.460999999999999
$ws.Range("D14").Value = -8.420899999999998
$ws.Range("A16").Value = -21.47400000000003
$ws.Range("B16").Value = 4.839300000000005
$ws.Range("D18").Value = -8.583899999999993
$ws.Range("A20").Value = -22.01870000000003
$ws.Range("B22").Value = 5.514500000000003
$ws.Range("D25").Value = -8.124899999999995
